# "Add files via upload" - append newly uploaded package names to the
# existing list of required R packages in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPackages = @("philentropy", "RcppML", "ggrepel", "concaveman", "ggforce")

# Find the first empty row after the existing data (column A) and append.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($i = 0; $i -lt $newPackages.Count; $i++) {
    $ws.Cells.Item($lastRow + 1 + $i, 1).Value = $newPackages[$i]
}

# Reproduce the scrolled position / selection left in the author's session.
$ws.Range("B54").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1

# The workbook was also set up for A4 portrait printing.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
